# Revert "Jeszcze jeden commit":
# Remove the three blank paragraphs and the "bbbbb...b" filler paragraph that were
# inserted after "Cwiczenia GIT", folding the trailing bookmark
# (_GoBack) back onto the "Cwiczenia GIT" paragraph.

$d = $word.ActiveDocument

# Drop the blank paragraphs sitting right after the first ("Cwiczenia GIT") paragraph.
while ($d.Paragraphs.Count -gt 2 -and $d.Paragraphs.Item(2).Range.Text.Trim().Length -eq 0) {
    $d.Paragraphs.Item(2).Range.Delete()
}

# Merge the following ("bbbbb...b") paragraph into the first one by deleting the
# paragraph mark that ends paragraph 1 — this keeps that paragraph's bookmark intact.
if ($d.Paragraphs.Count -gt 1) {
    $p1 = $d.Paragraphs.Item(1)
    $mark = $d.Range($p1.Range.End - 1, $p1.Range.End)
    $mark.Delete()
}

# Remove the filler run's text itself; the bookmarkStart/bookmarkEnd that followed it
# remain attached to the (now single) paragraph.
$d.Content.Find.Execute("bbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbbb", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null
